# Update cell values for rows 2-16 across columns B:F to reflect the new
# stock ticker lists, then remove the now-unused rows 17-26 so the sheet
# dimension shrinks from A1:F26 to A1:F16 (matching the target diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "NSE:ADSL"
$ws.Range("C2").Value = "NSE:ARTEMISMED"
$ws.Range("D2").Value = "NSE:LALPATHLAB"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "NSE:MGL"
$ws.Range("B3").Value = "NSE:AJOONI"
$ws.Range("C3").Value = "NSE:ARVINDFASN"
$ws.Range("D3").Value = "NSE:OFSS"
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = ""
$ws.Range("B4").Value = "NSE:ALANKIT"
$ws.Range("C4").Value = "NSE:FORTIS"
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("B5").Value = "NSE:BLUEDART"
$ws.Range("C5").Value = "NSE:IMAGICAA"
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("B6").Value = "NSE:CONSUMBEES"
$ws.Range("C6").Value = "NSE:INOXGREEN"
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = ""
$ws.Range("B7").Value = "NSE:DIGISPICE"
$ws.Range("C7").Value = "NSE:JINDWORLD"
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = ""
$ws.Range("B8").Value = "NSE:EBBETF0430"
$ws.Range("C8").Value = "NSE:KIRLPNU"
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("B9").Value = "NSE:HDFCLOWVOL"
$ws.Range("C9").Value = "NSE:KRITI"
$ws.Range("D9").Value = ""
$ws.Range("E9").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("B10").Value = "NSE:MGL"
$ws.Range("C10").Value = "NSE:MTNL"
$ws.Range("D10").Value = ""
$ws.Range("E10").Value = ""
$ws.Range("F10").Value = ""
$ws.Range("B11").Value = "NSE:MID150BEES"
$ws.Range("C11").Value = "NSE:NDTV"
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = ""
$ws.Range("B12").Value = "NSE:MOM100"
$ws.Range("C12").Value = "NSE:ORIENTPPR"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = ""
$ws.Range("B13").Value = "NSE:MONIFTY500"
$ws.Range("C13").Value = ""
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = ""
$ws.Range("F13").Value = ""
$ws.Range("B14").Value = "NSE:MONQ50"
$ws.Range("C14").Value = ""
$ws.Range("D14").Value = ""
$ws.Range("E14").Value = ""
$ws.Range("F14").Value = ""
$ws.Range("B15").Value = "NSE:RAMANEWS"
$ws.Range("C15").Value = ""
$ws.Range("D15").Value = ""
$ws.Range("E15").Value = ""
$ws.Range("F15").Value = ""
$ws.Range("B16").Value = "NSE:SAKUMA"
$ws.Range("C16").Value = ""
$ws.Range("D16").Value = ""
$ws.Range("E16").Value = ""
$ws.Range("F16").Value = ""

# Remove rows 17-26 entirely (their data moved up / was dropped).
$ws.Rows("17:26").Delete()
